$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 98, shifting rows 98:148 down to 99:149.
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with the new data entry.
$ws.Range("A98").Value = 8
$ws.Range("B98").Value = "Terminal La Palmera de La Serena"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = 44452
$ws.Range("E98").Value = 4
$ws.Range("F98").Value = 100112012
$ws.Range("G98").Value = "Espinaca"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 3200
$ws.Range("K98").Value = 400
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 450
$ws.Range("N98").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O98").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P98").Value = 900
$ws.Range("Q98").Value = 0.5
$ws.Range("R98").Value = "Hortaliza"
